$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.769.05'
$ws.Range("E2").Value = '  +2.66%  '

$ws.Range("D3").Value = '2.523.96'
$ws.Range("E3").Value = '  +1.20%  '

$ws.Range("D4").Value = '''1.00'

$ws.Range("D5").Value = '''592.25'
$ws.Range("E5").Value = '  +2.31%  '

$ws.Range("D6").Value = '''176.68'
$ws.Range("E6").Value = '  +3.82%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").Value = '''0.531'
$ws.Range("E8").Value = '  +2.60%  '

$ws.Range("D9").Value = '2.522.63'
$ws.Range("E9").Value = '  +1.18%  '

$ws.Range("E10").Value = '  +3.09%  '

$ws.Range("E11").Value = '  +3.07%  '

$ws.Range("D12").Value = '''5.15'
$ws.Range("E12").Value = '  +1.46%  '

$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").Value = '''26.83'
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").Value = '2.981.60'
$ws.Range("E15").Value = '  +1.43%  '

$ws.Range("E16").Value = '  +2.41%  '

$ws.Range("D17").Value = '67.545.38'
$ws.Range("E17").Value = '  +2.61%  '

$ws.Range("D18").Value = '2.515.04'
$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("E19").Value = '  +5.23%  '

$ws.Range("E20").Value = '  +2.57%  '

$ws.Range("D21").Value = '''363.21'
$ws.Range("E21").Value = '  +5.82%  '

$ws.Range("D22").Value = '''4.19'
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("E23").Value = '  +2.69%  '

$ws.Range("D24").Value = '''1.95'
$ws.Range("E24").Value = '  +1.04%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''71.06'
$ws.Range("E26").Value = '  +3.12%  '

$ws.Range("D27").Value = '''10.22'
$ws.Range("E27").Value = '  +3.92%  '

$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").Value = '2.653.70'
$ws.Range("E29").Value = '  +1.37%  '

$ws.Range("D30").Value = '0.0₃0988'
$ws.Range("E30").Value = '  +3.04%  '

$ws.Range("D31").Value = '''541.50'
$ws.Range("E31").Value = '  +3.69%  '

$ws.Range("D32").Value = '''8.29'
$ws.Range("E32").Value = '  +3.24%  '

$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("D34").Value = '''1.86'
$ws.Range("E34").Value = '  +2.99%  '

$ws.Range("D35").Value = '''0.129'
$ws.Range("E35").Value = '  -0.81%  '

$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  +1.24%  '

$ws.Range("D38").Value = '''156.34'
$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("D39").Value = '''18.83'
$ws.Range("E39").Value = '  +2.04%  '

$ws.Range("D40").Value = '''18.66'
$ws.Range("E40").Value = '  +1.89%  '

$ws.Range("E41").Value = '  +1.59%  '

$ws.Range("D42").Value = '''5.18'
$ws.Range("E42").Value = '  +2.74%  '

$ws.Range("D43").Value = '''1.79'
$ws.Range("E43").Value = '  +1.90%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.51'
$ws.Range("E45").Value = '  +3.53%  '

$ws.Range("D46").Value = '''0.561'
$ws.Range("E46").Value = '  +1.60%  '

$ws.Range("D47").Value = '''146.34'
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").Value = '''3.72'
$ws.Range("E48").Value = '  +2.03%  '

$ws.Range("D49").Value = '0.0₆0275'
$ws.Range("E49").Value = '  +2.30%  '

$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("D51").Value = '''0.0756'
$ws.Range("E51").Value = '  +0.85%  '
